$wb = $excel.ActiveWorkbook

# --- Update "Data" sheet: append 3 new daily observations ---
$dataSheet = $wb.Worksheets.Item("Data")

$dataSheet.Range("A441").Value = 45117
$dataSheet.Range("B441").Value = 1811.981

$dataSheet.Range("A442").Value = 45118
$dataSheet.Range("B442").Value = 1775.796

$dataSheet.Range("A443").Value = 45119
$dataSheet.Range("B443").Value = 1820.146

# Copy style of the previous data row (A440) to the new rows so date formatting stays consistent
$dataSheet.Range("A440").Copy()
$dataSheet.Range("A441:A443").PasteSpecial(-4122)

# --- Update "SeriesInfo" sheet: refresh metadata fields ---
# These cells hold plain text (not real dates), so force text entry and then
# strip the number-format Excel auto-applies for date-looking strings, to
# keep them as General-formatted text cells like the originals.
$seriesInfo = $wb.Worksheets.Item("SeriesInfo")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $seriesInfo.Range("B3") "2023-07-12"
Set-TextValue $seriesInfo.Range("B4") "2023-07-12"
Set-TextValue $seriesInfo.Range("B7") "2023-07-12"
Set-TextValue $seriesInfo.Range("B14") "2023-07-12 13:01:06-05"
